# Apply the latest cryptos list refresh (GitHub Actions scheduled update).
# For each changed row, update the Price (column D) and Volume(1h) (column E)
# cells. Column D values are written through a temporary "Text" number format
# so that numeric-looking strings (e.g. "0.999") are not auto-converted to
# numbers by Excel; ClearFormats() afterwards drops that temporary format so
# the cell's style index is left exactly as it was (no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '53.937.51'
$c.ClearFormats()
$ws.Range('E2').Value = '  -4.09%  '

# Row 3: Ethereum
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.248.74'
$c.ClearFormats()
$ws.Range('E3').Value = '  -5.08%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.11%  '

# Row 5: BNB
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '491.59'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.87%  '

# Row 6: Solana
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '127.07'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.42%  '

# Row 7: USDC
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E7').Value = '  +0.07%  '

# Row 8: XRP
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.528'
$c.ClearFormats()
$ws.Range('E8').Value = '  -2.84%  '

# Row 9: LidoStakedEther
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.271.86'
$c.ClearFormats()
$ws.Range('E9').Value = '  -4.37%  '

# Row 10: Dogecoin
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0931'
$c.ClearFormats()
$ws.Range('E10').Value = '  -5.63%  '

# Row 11: TRON
$ws.Range('E11').Value = '  -0.01%  '

# Row 12: Cardano
$ws.Range('E12').Value = '  -0.42%  '

# Row 13: Toncoin
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.67'
$c.ClearFormats()
$ws.Range('E13').Value = '  -4.48%  '

# Row 14: WrappedliquidstakedEther2.0
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.652.85'
$c.ClearFormats()
$ws.Range('E14').Value = '  -4.79%  '

# Row 15: Avalanche
$ws.Range('E15').Value = '  -1.05%  '

# Row 16: WrappedBTC
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '53.969.70'
$c.ClearFormats()
$ws.Range('E16').Value = '  -3.98%  '

# Row 17: ShibaInu
$ws.Range('E17').Value = '  -3.45%  '

# Row 18: WrappedEther
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.271.68'
$c.ClearFormats()
$ws.Range('E18').Value = '  -3.73%  '

# Row 19: Chainlink
$ws.Range('E19').Value = '  -2.12%  '

# Row 20: Polkadot
$ws.Range('E20').Value = '  +0.02%  '

# Row 21: BitcoinCash
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '298.13'
$c.ClearFormats()
$ws.Range('E21').Value = '  -3.74%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  -0.16%  '

# Row 23: Dai
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.12%  '

# Row 24: Litecoin
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '63.49'
$c.ClearFormats()
$ws.Range('E24').Value = '  -3.27%  '

# Row 25: Binance-PegBSC-USD
$ws.Range('E25').Value = '  +0.32%  '

# Row 26: Polygon
$ws.Range('E26').Value = '  +0.30%  '

# Row 27: WrappedeETH
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.368.04'
$c.ClearFormats()
$ws.Range('E27').Value = '  -4.50%  '

# Row 28: Kaspa
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.148'
$c.ClearFormats()
$ws.Range('E28').Value = '  +0.24%  '

# Row 29: InternetComputer(DFINITY)
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.13'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.93%  '

# Row 30: Monero
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '162.98'
$c.ClearFormats()
$ws.Range('E30').Value = '  -5.95%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  -2.99%  '

# Row 32: PEPE
$ws.Range('E32').Value = '  -4.59%  '

# Row 33: USDe
$ws.Range('E33').Value = '  -0.01%  '

# Row 34: Aptos
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.81'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.78%  '

# Row 35: FirstDigitalUSD
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.993'
$c.ClearFormats()
$ws.Range('E35').Value = '  -0.35%  '

# Row 36: Fetch.AI
$ws.Range('E36').Value = '  -0.83%  '

# Row 37: EthereumClassic
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '17.44'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.63%  '

# Row 38: ImmutableX
$ws.Range('E38').Value = '  +0.37%  '

# Row 39: SuiNetwork
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.839'
$c.ClearFormats()
$ws.Range('E39').Value = '  +1.58%  '

# Row 40: NEARProtocol
$ws.Range('E40').Value = '  -1.41%  '

# Row 41: OKB
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '35.52'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.19%  '

# Row 42: PolygonEcosystemToken
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.374'
$c.ClearFormats()
$ws.Range('E42').Value = '  +0.52%  '

# Row 43: Stacks
$ws.Range('E43').Value = '  +0.18%  '

# Row 44: Filecoin
$ws.Range('E44').Value = '  -1.24%  '

# Row 45: Aave
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '126.11'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.40%  '

# Row 46: RenderToken
$ws.Range('E46').Value = '  +1.28%  '

# Row 47: Stellar
$ws.Range('E47').Value = '  -0.98%  '

# Row 48: Bittensor
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '242.07'
$c.ClearFormats()
$ws.Range('E48').Value = '  +1.75%  '

# Row 49: Mantle
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.544'
$c.ClearFormats()
$ws.Range('E49').Value = '  -3.25%  '

# Row 50: Hedera
$ws.Range('E50').Value = '  -0.78%  '

# Row 51: VeChain
$ws.Range('E51').Value = '  -1.88%  '
